$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:C13").ClearContents()
$ws.Range("A11:B13").Select()
